$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Find the row that corresponds to "syntok" in column A and delete the entire row.
$found = $ws.Cells.Find("syntok", [Type]::Missing, [Type]::Missing, 1)
if ($found -ne $null) {
    $found.EntireRow.Delete()
}

$ws.Application.GoTo($ws.Range("D34"), $true)
